$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 13/14 by copying row 12 (ActorID) down, the same way the
# original sheet's existing rows were grown - this carries along styles and the
# "list" data validation so it naturally splits per-row like the source rows.
$ws.Rows("12").Copy() | Out-Null
$ws.Rows("13:14").Insert(-4121, 0) | Out-Null  # xlShiftDown, xlFormatFromLeftOrAbove

# Row 13: CamOffestPos
$ws.Range("A13").Value = "CamOffestPos"
$ws.Range("B13").Value = "string"
$ws.Range("C13").Value = $false
$ws.Range("D13").Value = $false
$ws.Range("E13").Value = $false
$ws.Range("F13").Value = $true
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = "Friend"
$ws.Range("J13").Value = "acctorid"

# Row 14: CamOffestRot
$ws.Range("A14").Value = "CamOffestRot"
$ws.Range("B14").Value = "string"
$ws.Range("C14").Value = $false
$ws.Range("D14").Value = $false
$ws.Range("E14").Value = $false
$ws.Range("F14").Value = $true
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = "Friend"
$ws.Range("J14").Value = "acctorid"

$ws.Range("A14").Select() | Out-Null

# The inserted rows' F cells don't carry the "TRUE,FALSE" list validation along
# automatically - apply it explicitly to the new rows.
$ws.Range("F13:F14").Validation.Add(3, 1, 1, '"TRUE,FALSE"') | Out-Null
